# Apply the text changes described in the diff:
#  - "SYSTEM Exibe os detalhes relativos àquela prestação de contas (dados básicos
#     da solicitação e documentos anexos); Exibe o histórico da tramitação da
#     prestação de contas."
#    becomes
#    "SYSTEM Exibe os detalhes relativos àquela prestação de contas (nome do
#     beneficiário, dados básicos da solicitação e documentos anexos); Exibe o
#     histórico da tramitação da prestação de contas."
#
#  - "Chefe Clica para analisar a prestação de contas."
#    becomes
#    "Chefe Verifica o histório da tramitação da prestação de contas e clica
#     para analisar a prestação de contas."
#
# These two strings appear (shared) in each of the four test-case blocks of the
# "Test Suite" sheet, at cells D19/D30/D41/D52 and B20/B31/B42/B53 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$oldDetails = "SYSTEM Exibe os detalhes relativos àquela prestação de contas (dados básicos da solicitação e documentos anexos); Exibe o histórico da tramitação da prestação de contas."
$newDetails = "SYSTEM Exibe os detalhes relativos àquela prestação de contas (nome do beneficiário, dados básicos da solicitação e documentos anexos); Exibe o histórico da tramitação da prestação de contas."

$oldClica = "Chefe Clica para analisar a prestação de contas."
$newClica = "Chefe Verifica o histório da tramitação da prestação de contas e clica para analisar a prestação de contas."

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -eq $oldDetails) {
            $cell.Value2 = $newDetails
        } elseif ($val -eq $oldClica) {
            $cell.Value2 = $newClica
        }
    }
}
